$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 2 (the sample/intern record row) with the new test values.
# Columns C..V correspond to FirstName, LastName, Dob, Gender, PhoneNumber,
# Address, City, State, Zip, Email, Languages, University, AdditionalCourses,
# SpecificSkills, Months, FromDate, ToDate, Interest, Location, Motivation.
$ws.Range("C2").Value = "o"
$ws.Range("D2").Value = "po"
$ws.Range("E2").Value = "po"
$ws.Range("F2").Value = "p"
$ws.Range("G2").Value = "op"
$ws.Range("H2").Value = "oo"
$ws.Range("I2").Value = "o"
$ws.Range("J2").Value = "p"
$ws.Range("K2").Value = "op"
$ws.Range("L2").Value = "o"
$ws.Range("M2").Value = "o"
$ws.Range("N2").Value = "o"
$ws.Range("O2").Value = "po"
$ws.Range("P2").Value = "o"
$ws.Range("Q2").Value = "op"
$ws.Range("R2").Value = "o"
$ws.Range("S2").Value = "o"
$ws.Range("T2").Value = "o"
$ws.Range("U2").Value = "o"
$ws.Range("V2").Value = "po"
